# Fix mojibake "Â±" -> "±" in the plus-minus statistics columns (B, C, D)
# for rows 2 through 17 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$badChar = [string][char]0x00C2 + [string][char]0x00B1   # "Â±"
$goodChar = [string][char]0x00B1                          # "±"

for ($row = 2; $row -le 17; $row++) {
    foreach ($col in @("B", "C", "D")) {
        $cell = $ws.Range("$col$row")
        $text = $cell.Text
        if ($text -ne $null -and $text.Contains($badChar)) {
            $cell.Value = $text.Replace($badChar, $goodChar)
        }
    }
}
